# Generating and pushing into server. Also deploys functions
#
# Adds four new parameter rows (bid / Colors / DefaultInTerminals /
# DefaultOutTerminals) to the "blockTemplate" sheet's template table, lets
# the sheet's helper formulas (columns D and I) extend down through the
# new rows, and leaves "blockTemplate" as the active/selected sheet & cell
# the way the author's Excel session ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("blockTemplate")

# Fill column A (the "name" column) first, then column B (the "type"
# column) - matches the order new entries were typed in the original
# editing session (and therefore the shared-string insertion order).
$ws.Range("A18").Value = "bid"
$ws.Range("A19").Value = "Colors"
$ws.Range("A20").Value = "DefaultInTerminals"
$ws.Range("A21").Value = "DefaultOutTerminals"

$ws.Range("B18").Value = "string"
$ws.Range("B19").Value = "object"
$ws.Range("B20").Value = "float"
$ws.Range("B21").Value = "float"

# Continue the D (JSON fragment) and I (trimmed name) helper formulas for
# the four new rows, same pattern used by every existing row.
for ($r = 18; $r -le 21; $r++) {
    $next = $r + 1
    $ws.Range("D" + $r).Formula = '="{"&CHAR(34)&$A$1&CHAR(34)&":"&CHAR(34)&A' + $r + '&CHAR(34)&","&CHAR(34)&$B$1&CHAR(34)&":"&CHAR(34)&B' + $r + '&CHAR(34)&","&CHAR(34)&$C$1&CHAR(34)&":"&CHAR(34)&C' + $r + '&CHAR(34)&"}"&IF(ISBLANK(A' + $next + '),"",",")'
    $ws.Range("I" + $r).Formula = "=TRIM(A" + $r + ")"
}

# Recalculate everything (E1's CONCAT(D:D) roll-up and D17's trailing-comma
# logic both depend on the freshly written rows below them).
$excel.CalculateFull()

# Leave blockTemplate as the active sheet/tab with B20 selected, as in the
# saved workbook.
$ws.Activate()
$ws.Range("B20").Select()
